$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D1" = 83.25813693486151
    "E1" = 71.67677128301931
    "F1" = 80.10421153793146
    "G1" = 49.77323929997574
    "D2" = 88.69064379680171
    "E2" = 83.5820895522388
    "F2" = 88.93408196366551
    "G2" = 57.15146087056657
    "D3" = 77.36234332043063
    "E3" = 82.7399820993484
    "F3" = 87.49168007674703
    "G3" = 56.59935659290447
    "D4" = 71.847952275506
    "E4" = 74.35897435897436
    "F4" = 78.90054256799399
    "G4" = 54.77615803240615
    "D5" = 85.6372827442285
    "E5" = 73.33333333333333
    "F5" = 87.97758469471823
    "G5" = 56.22302412499841
    "D6" = 76.01580307869507
    "E6" = 63.45373046756322
    "F6" = 66.82727256720689
    "G6" = 32.79735102212972
    "D7" = 85.23002923335929
    "E7" = 63.1578947368421
    "F7" = 82.6639771799408
    "G7" = 46.55771347202847
    "D8" = 81.28104102987085
    "E8" = 52.94117647058824
    "F8" = 62.68403347916787
    "G8" = 36.31741651201629
    "D9" = 100
    "E9" = 79.84698924526595
    "F9" = 85.35451977401129
    "G9" = 57.7634337727558
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
